# Re-order (permute) the data rows 149-170 on the "Artfynd" sheet.
# Each destination row ends up with the *entire* row content (all columns
# A:AY) that used to live in a different source row - i.e. the block of
# rows 149-170 gets shuffled into a new order while every row's own data
# (id, species, coordinates, etc.) travels with it unchanged.
#
# Mapping: new row R gets the content that currently sits in row Map[R].
# (rows 162 and 166 are fixed points - they keep their own content)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 149
$lastRow  = 170
$firstCol = "A"
$lastCol  = "AY"

# destination (array-relative, 1-based) -> source (array-relative, 1-based)
# relative row 1 == sheet row 149, relative row 22 == sheet row 170
$map = @{
    1  = 6
    2  = 13
    3  = 2
    4  = 10
    5  = 7
    6  = 16
    7  = 15
    8  = 20
    9  = 21
    10 = 19
    11 = 3
    12 = 11
    13 = 4
    14 = 14
    15 = 12
    16 = 8
    17 = 1
    18 = 18
    19 = 9
    20 = 17
    21 = 22
    22 = 5
}

# Text columns in this block must keep storing their (possibly
# numeric-/date-looking) content as literal text, e.g. "2" or
# "2023-08-26" - otherwise Excel auto-converts such values to a real
# number/date the moment they're written back via Value2. Number/boolean
# columns are left alone so they stay real numbers/booleans.
$textCols = @("C","D","F","G","H","I","J","K","N","P","T","U","V","W", `
              "Y","Z","AA","AB","AC","AF","AH","AI","AJ","AK","AO","AT", `
              "AW","AX","AY")
foreach ($tc in $textCols) {
    $ws.Range("$tc$firstRow`:$tc$lastRow").NumberFormat = "@"
}

$srcRange = $ws.Range("$firstCol$firstRow`:$lastCol$lastRow")
$orig = $srcRange.Value2

# $orig is a COM SAFEARRAY and is 1-based: [1..rows, 1..cols]
$rows = $orig.GetUpperBound(0)
$cols = $orig.GetUpperBound(1)

# a freshly allocated .NET array is 0-based, so build it with 0-based
# dimensions and offset the column/row indices accordingly when copying.
$new = New-Object 'object[,]' $rows, $cols

for ($r = 1; $r -le $rows; $r++) {
    $srcRow = $map[$r]
    for ($c = 1; $c -le $cols; $c++) {
        $new[($r - 1), ($c - 1)] = $orig[$srcRow, $c]
    }
}

$srcRange.Value2 = $new
